# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = -0.03276308475587503; "C2" = 3.590912207024844; "D2" = 47.40170072604877; "E2" = 6.884889303834068; "F2" = 6.953318621525447; "G2" = 51
    "B3" = 0.05364334786105588;  "C3" = 3.574536985163549; "D3" = 44.07793895074698; "E3" = 6.639121850873576; "F3" = 6.706306911391149; "G3" = 50
    "B4" = -0.04985138525762726; "C4" = 3.773358965575171; "D4" = 43.3640935270437;  "E4" = 6.585141875999613; "F4" = 6.653192858124505; "G4" = 49
    "B5" = 0.2063883000833151;   "C5" = 3.726029611715911; "D5" = 41.36251091033633; "E5" = 6.431369287355246; "F5" = 6.496080503111255; "G5" = 48
    "B6" = -0.03263097726318537; "C6" = 3.889169391173965; "D6" = 44.20905250212165; "E6" = 6.648988833057373; "F6" = 6.720790945475113; "G6" = 47
    "B7" = 0.2010956551540337;   "C7" = 3.91100830380399;  "D7" = 49.94279819706635; "E7" = 7.06702187608517;  "F7" = 7.142219545754433; "G7" = 46
    "B8" = -0.1796802820817339;  "C8" = 3.730114659890623; "D8" = 39.99999453259167; "E8" = 6.324554888100163; "F8" = 6.393439337599808; "G8" = 45
    "B9" = 0.01305864851553937;  "C9" = 3.680336663577857; "D9" = 42.37372025704332; "E9" = 6.509509985939289; "F9" = 6.584753686948571; "G9" = 44
    "B10" = -0.1400871512439874; "C10" = 4.022154036027911; "D10" = 46.29094780507954; "E10" = 6.803745130814318; "F10" = 6.882806227928441; "G10" = 43
    "B11" = 0.1042633263670164;  "C11" = 3.990286954105756; "D11" = 47.50789127154154; "E11" = 6.892596845278384; "F11" = 6.975348331320163; "G11" = 42
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
